$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 26
$ws.Cells.Item(5, 6).Value = 861
$ws.Cells.Item(6, 6).Value = 669
$ws.Cells.Item(7, 6).Value = 1239
$ws.Cells.Item(10, 6).Value = 704
$ws.Cells.Item(13, 6).Value = 365
$ws.Cells.Item(14, 6).Value = 726
$ws.Cells.Item(15, 6).Value = 965
$ws.Cells.Item(16, 6).Value = 10064
$ws.Cells.Item(17, 6).Value = 633
$ws.Cells.Item(19, 6).Value = 314
$ws.Cells.Item(22, 6).Value = 275
$ws.Cells.Item(29, 6).Value = 279
$ws.Cells.Item(30, 6).Value = 192
$ws.Cells.Item(31, 6).Value = 263
$ws.Cells.Item(32, 6).Value = 71
$ws.Cells.Item(33, 6).Value = 99
$ws.Cells.Item(35, 6).Value = 179
$ws.Cells.Item(37, 6).Value = 181
$ws.Cells.Item(38, 6).Value = 43

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 7).Value = "不可售"
$ws.Cells.Item(7, 6).Value = 138
$ws.Cells.Item(10, 6).Value = 243
$ws.Cells.Item(16, 6).Value = 292
$ws.Cells.Item(24, 6).Value = 8

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 822

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 822
$ws.Cells.Item(3, 6).Value = 26
$ws.Cells.Item(8, 3).Value = "广州·光启市七夕巧会·光夜only"
$ws.Cells.Item(8, 4).Value = "逸景路462号珠江国际纺织城d区6层 珠江时尚馆"
$ws.Cells.Item(8, 5).Value = "2024.08.10 10:00-08.10 21:00"
$ws.Cells.Item(8, 6).Value = 861
$ws.Cells.Item(8, 7).Value = 38
$ws.Cells.Item(8, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87920"
$ws.Cells.Item(8, 9).Value = "//i2.hdslb.com/bfs/openplatform/202406/eqyXGRCh1718954911688.png"
$ws.Cells.Item(9, 3).Value = "广州·次元喵喵动漫嘉年华"
$ws.Cells.Item(9, 4).Value = "东沙大道16号 广州健康方舟"
$ws.Cells.Item(9, 5).Value = "2024.08.10 10:00-08.10 21:00"
$ws.Cells.Item(9, 6).Value = 669
$ws.Cells.Item(9, 7).Value = 52
$ws.Cells.Item(9, 8).Value = "https://show.bilibili.com/platform/detail.html?id=88099"
$ws.Cells.Item(9, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/oilUtdD01718600424271.jpeg"
$ws.Cells.Item(10, 3).Value = "广州·火影忍者only"
$ws.Cells.Item(10, 4).Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws.Cells.Item(10, 5).Value = "2024.08.10 10:00-08.10 17:00"
$ws.Cells.Item(10, 6).Value = 1239
$ws.Cells.Item(10, 7).Value = 54
$ws.Cells.Item(10, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85704"
$ws.Cells.Item(10, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/lKOROXve1715763433389.jpeg"
$ws.Cells.Item(11, 3).Value = "广州·电子音乐仓库派对：RoaringLand咆哮大陆"
$ws.Cells.Item(11, 4).Value = "革新路124号太古仓码头4号仓 MAO Livehouse 广州(太古仓店)"
$ws.Cells.Item(11, 5).Value = "2024.08.10 23:00-08.11 04:00"
$ws.Cells.Item(11, 6).Value = 4
$ws.Cells.Item(11, 7).Value = 198
$ws.Cells.Item(11, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90061"
$ws.Cells.Item(11, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/HRdYR5SK1721892863746.png"
$ws.Cells.Item(13, 6).Value = 138
$ws.Cells.Item(15, 6).Value = 704
$ws.Cells.Item(17, 6).Value = 365
$ws.Cells.Item(19, 6).Value = 965
$ws.Cells.Item(20, 6).Value = 10065
$ws.Cells.Item(21, 6).Value = 243
$ws.Cells.Item(22, 6).Value = 633
$ws.Cells.Item(25, 6).Value = 275
$ws.Cells.Item(36, 6).Value = 279
$ws.Cells.Item(37, 6).Value = 192
$ws.Cells.Item(38, 6).Value = 263
$ws.Cells.Item(39, 6).Value = 71
$ws.Cells.Item(40, 6).Value = 99
$ws.Cells.Item(43, 6).Value = 179
$ws.Cells.Item(47, 6).Value = 181
$ws.Cells.Item(50, 6).Value = 8

Write-Host "Edits applied successfully"